$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recomputed loading_percent results for the "380 kV" case (rows 2-25,
# i.e. time steps 0-23). Only the line-loading columns that actually
# change are touched (B, C, D, E, G, H, I, J, O); F, K, L, M, N stay 0.

$col = New-Object 'object[,]' 24,1
$col[0,0] = 16.19292459086031
$col[1,0] = 15.432680365463
$col[2,0] = 14.9458988786319
$col[3,0] = 14.74272291016459
$col[4,0] = 14.70870231809528
$col[5,0] = 14.94317793302778
$col[6,0] = 15.93506117513978
$col[7,0] = 17.71361224664487
$col[8,0] = 18.91014639794101
$col[9,0] = 19.42920841481017
$col[10,0] = 19.62204577407954
$col[11,0] = 19.58068151062172
$col[12,0] = 19.44514827422049
$col[13,0] = 19.36164339465563
$col[14,0] = 18.87570839600928
$col[15,0] = 18.57106748443991
$col[16,0] = 18.39347390736037
$col[17,0] = 18.33293935146496
$col[18,0] = 18.60374324238435
$col[19,0] = 19.4850592697996
$col[20,0] = 20.03933558808188
$col[21,0] = 19.74552033521141
$col[22,0] = 18.5889781610148
$col[23,0] = 17.25128556918714
$ws.Range("B2:B25").Value = $col

$col = New-Object 'object[,]' 24,1
$col[0,0] = 11.10604247870716
$col[1,0] = 10.45713507579169
$col[2,0] = 10.03657523240876
$col[3,0] = 9.859713467324447
$col[4,0] = 9.830017578427753
$col[5,0] = 10.0342120759383
$col[6,0] = 10.88693629131274
$col[7,0] = 12.38077941009831
$col[8,0] = 13.36689941352887
$col[9,0] = 13.79098092761027
$col[10,0] = 13.94802703178554
$col[11,0] = 13.9143622283031
$col[12,0] = 13.803972330429
$col[13,0] = 13.73589327046192
$col[14,0] = 13.33868963914241
$col[15,0] = 13.08872426015816
$col[16,0] = 12.94264455596891
$col[17,0] = 12.89278909977518
$col[18,0] = 13.11557235299481
$col[19,0] = 13.83649283882946
$col[20,0] = 14.28699190517767
$col[21,0] = 14.04844780069682
$col[22,0] = 13.10344171471739
$col[23,0] = 11.99607291901273
$ws.Range("C2:C25").Value = $col

$col = New-Object 'object[,]' 24,1
$col[0,0] = 14.8455803421141
$col[1,0] = 14.7792772947598
$col[2,0] = 14.74186032583804
$col[3,0] = 14.72745250121655
$col[4,0] = 14.72511114131523
$col[5,0] = 14.74166260150619
$col[6,0] = 14.82204290276501
$col[7,0] = 15.00523214745036
$col[8,0] = 15.1545746549742
$col[9,0] = 15.2255206544989
$col[10,0] = 15.25280083542013
$col[11,0] = 15.24690742549625
$col[12,0] = 15.22775681210081
$col[13,0] = 15.21607992967335
$col[14,0] = 15.14999717855406
$col[15,0] = 15.11021556477444
$col[16,0] = 15.08761863156008
$col[17,0] = 15.08001708800867
$col[18,0] = 15.11442108117771
$col[19,0] = 15.23337071217189
$col[20,0] = 15.31351804096206
$col[21,0] = 15.27052789925877
$col[22,0] = 15.11251891221399
$col[23,0] = 14.95301959731817
$ws.Range("D2:D25").Value = $col

$col = New-Object 'object[,]' 24,1
$col[0,0] = 16.24748311274231
$col[1,0] = 16.18117960467449
$col[2,0] = 16.14411074116382
$col[3,0] = 16.12993185548984
$col[4,0] = 16.12763375391441
$col[5,0] = 16.14391575251711
$col[6,0] = 16.22387291092954
$col[7,0] = 16.40901713336517
$col[8,0] = 16.56151696806935
$col[9,0] = 16.63427635601287
$col[10,0] = 16.66229757798371
$col[11,0] = 16.65624215338942
$col[12,0] = 16.63657238986948
$col[13,0] = 16.62458459309157
$col[14,0] = 16.55682858336014
$col[15,0] = 16.51611718153361
$col[16,0] = 16.49302095202934
$col[17,0] = 16.48525644607966
$col[18,0] = 16.52041800087036
$col[19,0] = 16.64233730134185
$col[20,0] = 16.72474226724418
$col[21,0] = 16.68051817929728
$col[22,0] = 16.51847263566404
$col[23,0] = 16.35597486941472
$ws.Range("E2:E25").Value = $col

$col = New-Object 'object[,]' 24,1
$col[0,0] = 27.97909451974359
$col[1,0] = 28.09123420176989
$col[2,0] = 28.17356870418759
$col[3,0] = 28.21048127797306
$col[4,0] = 28.21681266984735
$col[5,0] = 28.17405295296095
$col[6,0] = 28.01494546075791
$col[7,0] = 27.811109783795
$col[8,0] = 27.72881607412194
$col[9,0] = 27.7063056236163
$col[10,0] = 27.69994675056454
$col[11,0] = 27.7012196680724
$col[12,0] = 27.70573897878478
$col[13,0] = 27.70878969971163
$col[14,0] = 27.73058921520664
$col[15,0] = 27.74779969035685
$col[16,0] = 27.75910256995059
$col[17,0] = 27.76317002116985
$col[18,0] = 27.7458221637633
$col[19,0] = 27.70435264203561
$col[20,0] = 27.68987812674984
$col[21,0] = 27.69644236209484
$col[22,0] = 27.74671181804856
$col[23,0] = 27.85450668910556
$ws.Range("G2:G25").Value = $col

$col = New-Object 'object[,]' 24,1
$col[0,0] = 14.08430870042998
$col[1,0] = 14.15234660269855
$col[2,0] = 14.19728410152631
$col[3,0] = 14.21639027623882
$col[4,0] = 14.21961074619763
$col[5,0] = 14.19753856184676
$col[6,0] = 14.10711110913554
$col[7,0] = 13.95493120571188
$col[8,0] = 13.85853752454554
$col[9,0] = 13.81805050374204
$col[10,0] = 13.80320432953743
$col[11,0] = 13.80638010243372
$col[12,0] = 13.81681936190753
$col[13,0] = 13.82327697451889
$col[14,0] = 13.86125127042037
$col[15,0] = 13.88540987048092
$col[16,0] = 13.89962168575926
$col[17,0] = 13.90448785476236
$col[18,0] = 13.88280538427261
$col[19,0] = 13.81373991078994
$col[20,0] = 13.7714316599566
$col[21,0] = 13.79375282215762
$col[22,0] = 13.88398186735103
$col[23,0] = 13.99339927526547
$ws.Range("H2:H25").Value = $col

$col = New-Object 'object[,]' 24,1
$col[0,0] = 19.01267820699817
$col[1,0] = 19.17856015529042
$col[2,0] = 19.28579696678489
$col[3,0] = 19.33085247603309
$col[4,0] = 19.338415835678
$col[5,0] = 19.28639911002504
$col[6,0] = 19.06875769444949
$col[7,0] = 18.68461377918781
$col[8,0] = 18.42828057701817
$col[9,0] = 18.31727133450416
$col[10,0] = 18.27603889190185
$col[11,0] = 18.28488327125296
$col[12,0] = 18.31386299947122
$col[13,0] = 18.33171865304028
$col[14,0] = 18.43564788068102
$col[15,0] = 18.50083838363872
$col[16,0] = 18.53886117825355
$col[17,0] = 18.55182556768603
$col[18,0] = 18.49384421186721
$col[19,0] = 18.30532912544147
$col[20,0] = 18.18681177917308
$col[21,0] = 18.24963792558402
$col[22,0] = 18.49700458257427
$col[23,0] = 18.78398001240804
$ws.Range("I2:I25").Value = $col

$col = New-Object 'object[,]' 24,1
$col[0,0] = 9.260850407519932
$col[1,0] = 9.267779451442314
$col[2,0] = 9.273461692722238
$col[3,0] = 9.276136457802862
$col[4,0] = 9.276602299175654
$col[5,0] = 9.273496310958683
$col[6,0] = 9.2629433365572
$col[7,0] = 9.253567512635106
$col[8,0] = 9.253561731617419
$col[9,0] = 9.255048044525964
$col[10,0] = 9.2558243175234
$col[11,0] = 9.255647652251922
$col[12,0] = 9.255107635542663
$col[13,0] = 9.254804634107902
$col[14,0] = 9.253494503990256
$col[15,0] = 9.253071691386628
$col[16,0] = 9.252968693581685
$col[17,0] = 9.252957917287134
$col[18,0] = 9.253102195606163
$col[19,0] = 9.255260464720074
$col[20,0] = 9.257914719062025
$col[21,0] = 9.256384542389037
$col[22,0] = 9.253087968253459
$col[23,0] = 9.253567512635106
$ws.Range("J2:J25").Value = $col

$col = New-Object 'object[,]' 24,1
$col[0,0] = 21.33706928720452
$col[1,0] = 21.44803872421981
$col[2,0] = 21.5227719296964
$col[3,0] = 21.55487758731322
$col[4,0] = 21.56030821284759
$col[5,0] = 21.52319824219246
$col[6,0] = 21.37395760414579
$col[7,0] = 21.13398833685594
$col[8,0] = 20.99026917179188
$col[9,0] = 20.93206123393574
$col[10,0] = 20.91105847665222
$col[11,0] = 20.91553544283655
$col[12,0] = 20.93031244834693
$col[13,0] = 20.93949938276912
$col[14,0] = 20.99421819258124
$col[15,0] = 21.02962885419903
$col[16,0] = 21.05067058125577
$col[17,0] = 21.05791053426872
$col[18,0] = 21.02578946582874
$col[19,0] = 20.92594381196603
$col[20,0] = 20.86675059454625
$col[21,0] = 20.89778581097265
$col[22,0] = 21.02752312430047
$col[23,0] = 21.19321717036284
$ws.Range("O2:O25").Value = $col

Write-Output "Updated loading_percent values for rows 2-25 (columns B,C,D,E,G,H,I,J,O)."
